# Update crypto price/volume data per upstream refresh (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.600.89"
$ws.Range("E2").Value = "  +2.32%  "
$ws.Range("D3").Value = "1.673.49"
$ws.Range("E3").Value = "  +2.70%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.78"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.58%  "
$ws.Range("E6").Value = "  +2.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.63"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.95%  "
$ws.Range("E9").Value = "  +2.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0645"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +6.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0905"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "1.913.93"
$ws.Range("E12").Value = "  +2.61%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.28"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +12.75%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.672.09"
$ws.Range("E14").Value = "  +2.60%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.614"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +9.44%  "
$ws.Range("E16").Value = "  +4.46%  "
$ws.Range("D17").Value = "30.614.09"
$ws.Range("E17").Value = "  +2.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.38"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.36"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").Value = "0.0₃0724"
$ws.Range("E20").Value = "  +3.31%  "
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("E22").Value = "  +3.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.00"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.17"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.47"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.89"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.68%  "
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("E31").Value = "  +3.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.46"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.18%  "
$ws.Range("E33").Value = "  +4.15%  "
$ws.Range("D34").Value = "1.485.13"
$ws.Range("E34").Value = "  +4.35%  "
$ws.Range("E35").Value = "  +7.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "84.89"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +12.90%  "
$ws.Range("E37").Value = "  -0.57%  "
$ws.Range("E38").Value = "  +9.10%  "
$ws.Range("E39").Value = "  +5.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.69"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.98%  "
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.841"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("E43").Value = "  +1.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.98"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "51.57"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.51"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.17%  "
$ws.Range("D49").Value = "1.807.24"
$ws.Range("E49").Value = "  +1.94%  "
$ws.Range("E50").Value = "  +4.88%  "
$ws.Range("E51").Value = "  -0.34%  "
